$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $result = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $result) {
        Write-Host "WARNING: replace failed for:" $findText
    }
}

function Find-Paragraph($matchText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $matchText) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphPlainText($para, $newText) {
    # Remove all text/formatting from the paragraph's run(s) (but keep the
    # paragraph mark / pPr intact), then insert fresh, unformatted text.
    # This avoids carrying over run-level rPr (e.g. italics) on the new text.
    $start = $para.Range.Start
    $end = $para.Range.End - 1
    if ($end -gt $start) {
        $r = $d.Range($start, $end)
        $r.Delete()
    }
    $insertPoint = $d.Range($start, $start)
    $insertPoint.InsertAfter($newText)
}

# ---------------------------------------------------------------------
# 1. Date change
# ---------------------------------------------------------------------
Replace-Text "2019-09-11" "2019-09-12"

# ---------------------------------------------------------------------
# 2. Ecosystem monitoring paragraph edits
# ---------------------------------------------------------------------
Replace-Text "response to changes in conditions and the accessibility" "response to changes in abiotic conditions, and the accessibility"
Replace-Text "collected from various studies, much of" "collected from various programs, much of"

# ---------------------------------------------------------------------
# 3. Key Largo, Florida -> Key Largo
# ---------------------------------------------------------------------
Replace-Text "Key Largo, Florida is the northernmost" "Key Largo is the northernmost"

# ---------------------------------------------------------------------
# 4. "The waters surrounding Key Largo..." paragraph edits
# ---------------------------------------------------------------------
Replace-Text "with one another through the use of water-mediating" "with one another by means of water-mediating"
Replace-Text "the Atlantic Ocean. The Florida Bay is a relatively small" "the Atlantic Ocean. Florida Bay is a relatively small"
Replace-Text "The bayside is a smaller, relatively shallow body, and enclosed body of water with a generally dynamic range" "The bayside is a relatively shallow enclosed body of water with a dynamic range"

# ---------------------------------------------------------------------
# 5. Replace the placeholder paragraph with real paragraph text
# ---------------------------------------------------------------------
Replace-Text "#######Add one more paragraph about Largo specific monitoring/this data to transition to the description of data and objectives/hypothesis.###" "The health of these three major aquatic ecosystems is dependant on the stability of the abiotic conditions of associated waters. Though each of these aquatic habitats are closely associted to the island of Key Largo, chemical and geographical conditions differ greatly between systems and thus support considerably different community structures. Adequate monitoring of these coastal water systems is imperative to the continued envrionmental and economic health of the island. Citizen science data collection enables reseachers to collect information at a substantially larger scale than would be possible alone. The analysis of these data will be used to establish a baseline health condition for various water bodies associated with Key Largo; and will produce data visuals to enrich future citizen science programs."

# ---------------------------------------------------------------------
# 6. Objective 4 title change
# ---------------------------------------------------------------------
Replace-Text "Objective 4: Giving Back to Marinelab" "Objective 4: Promoting Future Citizen Science"

# ---------------------------------------------------------------------
# 7. "Methods and Results" Heading1 -> "Methods" Heading2, bookmark rename
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("methods-and-results")
$bmRange = $bm.Range
$bm.Delete()
$bmRange.Text = "Methods"
$d.Bookmarks.Add("methods", $bmRange)

$methodsPara = Find-Paragraph "Methods*"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Methods`r") {
        $p.Style = "Heading 2"
    }
}

# ---------------------------------------------------------------------
# 8. Replace italic placeholder text under "Methods" with real content
# ---------------------------------------------------------------------
$p1 = Find-Paragraph "*In most research papers*"
Set-ParagraphPlainText $p1 "All data was collected by students or instructors of the Marinelab envrionmental education program. All collectors are required to complete a one hour training program to familarize themselves with the testing materials before they are permitted to collect data. Data is collected in the field at various sampling locations. Sample location is dictated by the needs of the program, accessibility of location, and weather conditions. On site, a grab sample of water is collected and brought onto the boat. Abiotic conditions are measured immediately follwoing collection and scored onto a paper record sheet."

# ---------------------------------------------------------------------
# 9 & 10. Insert two new BodyText paragraphs before "Data aquisition"
# heading (equivalent to right after the "All data was collected..."
# paragraph). Inserting via InsertBefore on the heading paragraph (rather
# than InsertAfter on the prior paragraph) avoids the new text getting
# swallowed into the "data-aquisition" bookmark, which starts exactly at
# that boundary.
# ---------------------------------------------------------------------
$dataAcqPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Data aquisition`r") {
        $dataAcqPara = $p
    }
}
$insertPoint = $d.Range($dataAcqPara.Range.Start, $dataAcqPara.Range.Start)
$insertPoint.InsertBefore("####Insert Parameter Specifics#######`rAll measures are taken using semiquantative testing equipment or with a YSI Sonde. All data sheets are collected and entered into the master database upon return to shore.`r")

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "####Insert Parameter Specifics#######`r") {
        $p.Style = "Body Text"
    }
    if ($p.Range.Text -eq "All measures are taken using semiquantative testing equipment or with a YSI Sonde. All data sheets are collected and entered into the master database upon return to shore.`r") {
        $p.Style = "Body Text"
    }
}

# ---------------------------------------------------------------------
# 11. Replace italic placeholder text under "Data aquisition" with real content
# ---------------------------------------------------------------------
$p2 = Find-Paragraph "*As applicable, explain where and how*"
Set-ParagraphPlainText $p2 "These data were aquired from the Marinelab citizen science master database. The database is privately maintained by the Marine Resources Development Foundation and is not readily available online. Data is primairly used as a tool to enrich scientific education in the Marinelab program, however data may be distrubuted for analysis at the discretion of the director of the program. This analysis was conducted under the permission of the current Marinelab director Sarah Egner, and the resulting analysis will be used for the enrichment of future programs."

# ---------------------------------------------------------------------
# 12. Insert new BodyText paragraph "##Results" before "Discussion" heading
# ---------------------------------------------------------------------
$discussionPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Discussion`r" -and $p.Style.NameLocal -eq "Heading 1") {
        $discussionPara = $p
    }
}
$insertPoint2 = $d.Range($discussionPara.Range.Start, $discussionPara.Range.Start)
$insertPoint2.InsertBefore("##Results`r")

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "##Results`r") {
        $p.Style = "Body Text"
    }
}

# ---------------------------------------------------------------------
# 13. "Discussion" Heading1 -> Heading2
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Discussion`r" -and $p.Style.NameLocal -eq "Heading 1") {
        $p.Style = "Heading 2"
    }
}

Write-Host "Edit complete"
